$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values are numeric-looking strings stored as text (matching the
# original inline-string cells). Force text number format before assigning
# so Excel does not silently convert them into floating point numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "246.77"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "21.99"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.405"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05795"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.334"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8073"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9602"
$ws.Range("E9").Value = "8FTXTokenFTT"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1429"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07520"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03217"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03019"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.138"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09407"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001589"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04801"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0005901"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.005616"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0009940"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.752"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.239"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3233"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1260"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03896"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006371"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1075"
$ws.Range("E43").Value = "42CEJICEJIBestin24h"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1467"
